$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.091.45"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.366.11"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.69"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.37%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.365.55"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.58"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.77%  "

$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("E12").Value = "  +3.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.939.42"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.41%  "

$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.374.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.23"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.205.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.61%  "

$ws.Range("E21").Value = "  +3.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "379.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +7.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.567"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.500.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.54"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +10.47%  "

$ws.Range("E28").Value = "  +13.36%  "

$ws.Range("E29").Value = "  +7.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.33%  "

$ws.Range("E32").Value = "  +3.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.13"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.393.70"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("E36").Value = "  +3.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.59%  "

$ws.Range("E38").Value = "  +2.77%  "

$ws.Range("E39").Value = "  +4.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "160.44"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("E41").Value = "  +2.91%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +10.20%  "

$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("E46").Value = "  +2.15%  "

$ws.Range("E47").Value = "  +7.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.03"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.92%  "

$ws.Range("E49").Value = "  +3.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.323.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.46%  "
